$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Activate()
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = "2016-08-24 02:16:19"
$ws.Range("G5").Value = "2016-08-24 02:16:19"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Activate()
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "mt"
$ws.Range("E5").Value = "mt"
$ws.Range("H2").Value = "2016-08-24 02:16:15"
$ws.Range("H5").Value = "2016-08-24 02:16:15"
$ws.Range("K2").Value = "2016-08-24 02:16:32"
$ws.Range("K5").Value = "2016-08-24 02:16:32"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Activate()
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "mt"
$ws.Range("E5").Value = "mt"
$ws.Range("H2").Value = "2016-08-24 02:16:19"
$ws.Range("H5").Value = "2016-08-24 02:16:19"
$ws.Range("K2").Value = "2016-08-24 02:16:39"
$ws.Range("K5").Value = "2016-08-24 02:16:39"

# Restore original active sheet (Overview was first sheet / active in the source file)
$wsOverview.Activate()
